# Iteration v0.9 -- Porter Stemmer
# Adds a new "0.9a"/"0.9b" iteration column-pair to the Results sheet
# (copy of the 0.7 column, now holding the new run's numbers) and records
# the corresponding pre-processing/feature-engineering/modeling steps
# (Wordninja, PorterStemmer) on the Steps sheet.

$wb = $excel.ActiveWorkbook
$wsResults = $wb.Worksheets.Item("Results")
$wsSteps   = $wb.Worksheets.Item("Steps")

# --- Results sheet: header row for the third mini-table (row 16) ---
$wsResults.Range("H16").Value = 0.8
$wsResults.Range("J16").Value = "0.9a"
$wsResults.Range("L16").Value = "0.9b"

# --- Results sheet: Accuracy row (18) ---
$wsResults.Range("H18").Value = 0.908207070707071
$wsResults.Range("J18").Value = 0.913636363636364
$wsResults.Range("L18").Value = 0.971085858585859

# --- Results sheet: FPR row (19) ---
$wsResults.Range("H19").Value = 0.0577492596248766
$wsResults.Range("J19").Value = 0.052319842053307
$wsResults.Range("L19").Value = 0.0498519249753208

# --- Results sheet: F1 row (20), Train/Test pairs ---
$wsResults.Range("H20").Value = 0.84004400440044
$wsResults.Range("I20").Value = 0.876056181061124
$wsResults.Range("J20").Value = 0.848806366047745
$wsResults.Range("L20").Value = 0.943858788918853
$wsResults.Range("M20").Value = 0.855143054009646

# --- Steps sheet: new rows describing the 0.8 / 0.9a / 0.9b iterations ---
$wsSteps.Range("A29").Value = 0.8
$wsSteps.Range("B29").Value = "Wordninja"

$wsSteps.Range("A31").Value = "0.9a"
$wsSteps.Range("B31").Value = "PorterStemmer"

$wsSteps.Range("A32").Value = "0.9b"
$wsSteps.Range("B32").Value = "PorterStemmer"
$wsSteps.Range("D32").Value = "Logistic Regression"

# --- Restore the selections left by the editor ---
# (Steps sheet selection is set first; Results stays the active/visible tab.)
[void]$wsSteps.Activate()
$stepsSelection = $excel.Union($wsSteps.Range("D32"), $wsSteps.Range("L16"))
[void]$stepsSelection.Select()

[void]$wsResults.Activate()
[void]$wsResults.Range("L16").Select()
